$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10Nov2019")

# ---------------------------------------------------------------
# 0. Capture values we need to relocate before the row-insert
#    (row 26 and row 27 are untouched by inserting at row 28).
# ---------------------------------------------------------------
$a27 = $ws.Range("A27").Value2()
$b27 = $ws.Range("B27").Value2()
$c27 = $ws.Range("C27").Value2()
$k27 = $ws.Range("K27").Value2()
$l27 = $ws.Range("L27").Value2()
$m27 = $ws.Range("M27").Value2()
$n27 = $ws.Range("N27").Value2()
$o27 = $ws.Range("O27").Value2()

# ---------------------------------------------------------------
# 1. Insert 3 new blank rows right after row 27 (new rows 28-30).
#    Excel clones row 27's cell formatting into the new rows, and
#    shifts the old rows 28+ down to 31+.
# ---------------------------------------------------------------
$ws.Rows("28:30").Insert()

# ---------------------------------------------------------------
# 2. Populate new row 30 with what used to be row 27's content.
#    (Formatting already matches thanks to the row-insert clone.)
# ---------------------------------------------------------------
$ws.Range("A30").Value = $a27
$ws.Range("B30").Value = $b27
$ws.Range("C30").Value = $c27
$ws.Range("K30").Value = $k27
$ws.Range("L30").Value = $l27
$ws.Range("M30").Value = $m27
$ws.Range("N30").Value = $n27
$ws.Range("O30").Value = $o27

# ---------------------------------------------------------------
# 3. Fix up new row 29 so that it matches the row-26-before-edit
#    look (F:H/K:M/N:O = percent style "41", Q:S = style "6").
# ---------------------------------------------------------------
$ws.Range("Q26:S26").Copy()
$ws.Range("Q29:S29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F26:H26").Copy()
$ws.Range("N29:O29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 4. Fix up new row 28 (K:M -> style 45, N:O -> style 46), using
#    row 25 (unchanged) as the style source.
# ---------------------------------------------------------------
$ws.Range("K25:O25").Copy()
$ws.Range("K28:O28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 5. Fix up row 26 (K:M -> style 45, N:O -> style 46, Q:S -> 41).
# ---------------------------------------------------------------
$ws.Range("K25:O25").Copy()
$ws.Range("K26:O26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F26:H26").Copy()
$ws.Range("Q26:S26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 6. Rebuild row 27: A:C take the "header row" look (style 24,
#    copied from row 24), C27 now references string 139 and a new
#    D27 cell (style 57, copied from D20) is added; K:O are
#    cleared but keep the style-45/46 percent look.
# ---------------------------------------------------------------
$ws.Range("A24:C24").Copy()
$ws.Range("A27:C27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D20").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A27").Value = "yes"
$ws.Range("B27").Value = "yes "
$ws.Range("C27").Value = "bert12->freeze->layernorm"
$ws.Range("D27").Value = "add 20 topology"

$ws.Range("K25:O25").Copy()
$ws.Range("K27:O27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("K27").Value = ""
$ws.Range("L27").Value = ""
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""
$ws.Range("O27").Value = ""

# ---------------------------------------------------------------
# 7. Add the new B34 cell (shared string "`"), re-using the style
#    of the already-present style-28 text cells on that row band.
# ---------------------------------------------------------------
$ws.Range("B36").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B34").Value = "``"

# ---------------------------------------------------------------
# 8. Restore the selection shown in the target workbook.
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("F24:H24").Select()
